$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("ServerDatabaseConfig")
$wsSet1   = $wb.Worksheets.Item("ScriptSet1")
$wsSet2   = $wb.Worksheets.Item("ScriptSet2")

# --- ServerDatabaseConfig: QueryTimeOut value 60 -> 30 ---
$wsConfig.Range("D2").Value = 30

# --- ScriptSet1: update NumberExec (F) values, rename two scripts, drop the two extra rows ---
$wsSet1.Range("F2").Value = 2
$wsSet1.Range("F3").Value = 1
$wsSet1.Range("F4").Value = 1
$wsSet1.Range("F5").Value = 1
$wsSet1.Range("D6").Value = "GetTodaysQueries.sql"
$wsSet1.Range("D7").Value = "BadQuery.sql"
$wsSet1.Rows("8:9").Delete()

# --- ScriptSet2: same NumberExec (F) + script-name updates, drop the two extra rows ---
$wsSet2.Range("F2").Value = 2
$wsSet2.Range("F3").Value = 1
$wsSet2.Range("F4").Value = 1
$wsSet2.Range("F5").Value = 1
$wsSet2.Range("D6").Value = "GetTodaysQueries.sql"
$wsSet2.Range("D7").Value = "BadQuery.sql"
$wsSet2.Rows("8:9").Delete()

# --- Selections / active sheet, matching the saved workbook UI state ---
$wsConfig.Range("D6").Select()
$wsSet2.Range("F2:F7").Select()

$wsSet1.Activate()
$wsSet1.Range("F13").Select()
